$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "What If? 2 English Hardcover" entry currently sitting in row 4 gets
# pushed down to the new last row (50), and a brand new wishlist entry
# ("Root Boardgame") takes its place in row 4.
# ---------------------------------------------------------------------------

# 1) Append the old row-4 data as the new row 50 (values + hyperlink style).
$ws.Range("A50").Value = "What If? 2 English Hardcover"
$ws.Range("B50").Value = "https://m.media-amazon.com/images/P/0593542908.01._SCLZZZZZZZ_SX500_.jpg"
$ws.Range("C50").Value = "https://www.amazon.de/-/en/Randall-Munroe/dp/147368062X/ref=sr_1_1?keywords=what+if+2+english&qid=1688464961&sprefix=what+if+2+e%2Caps%2C92&sr=8-1"
$ws.Range("D50").Value = "19.15 EUR"
$ws.Range("B50").Style = "Hyperlink"
$ws.Range("C50").Style = "Hyperlink"

# 2) Overwrite row 4 with the new "Root Boardgame" entry (plain style, no
#    hyperlink formatting yet). Written in Name/Link/Price/Image order so
#    new shared-string ids come out the same as the authored workbook.
$ws.Range("A4").Value = "Root Boardgame"
$ws.Range("C4").Value = "https://www.amazon.de/-/en/Leder-Games-LED01000-Root-Game/dp/B07F454YF3/ref=sr_1_2_mod_primary_new?dib=eyJ2IjoiMSJ9.6CDHwDz4ym78N7qETcYkjfoSY2mwJydyxT12aGMYUEbjeDX5EJUI2V7Nt04l1GQaXVxQVvc2WVVSqL5ZOlTzoxNHIRxNdbPa8SE0HpjF7vOM9Ij8WNlHurk5SVGm0IMoQwaDspMthHEbHHq5WrP4jfivUJll5JLDW4-ih0KI36X0IVBiIxEf24HpQKa_-neD.tSUwCbfvn4n4eabrwxs2YrOsGUGmbZuewzhphCWX8TQ&dib_tag=se&keywords=Leder+Games+Root&linkCode=gg3&qid=1750840429&sbo=RZvfv%2F%2FHxDF%2BO5021pAnSA%3D%3D&sr=8-2"
$ws.Range("D4").Value = "73 EUR"
$ws.Range("B4").Value = "https://m.media-amazon.com/images/I/91ezFG-gQ6L._AC_SX679_.jpg"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Style = "Normal"

# ---------------------------------------------------------------------------
# Hyperlinks: the engine only exposes a whole-sheet Hyperlinks.Delete(), so
# rebuild the full set (same targets, same order -> same rIds) with the
# "What If?" link now anchored on B50 instead of B4.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.wog.ch/nas/cover_large/p5/p5_finalfantasy16steelbook.jpg")
$ws.Hyperlinks.Add($ws.Range("B50"), "https://m.media-amazon.com/images/P/0593542908.01._SCLZZZZZZZ_SX500_.jpg")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.wog.ch/nas/cover_xl/pc/pc_citiesskylines2d1.jpg")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://store.steampowered.com/app/2050650/Resident_Evil_4/")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://cdn.akamai.steamstatic.com/steam/apps/2050650/header.jpg?t=1687479023")
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.loveramics.com/collections/diner-mugs/products/bond-250ml-starsky-mug-granite-carmel-gunpowder")
$ws.Hyperlinks.Add($ws.Range("C16"), "https://www.digitec.ch/en/s1/product/8bitdo-ultimate-controller-with-charging-station-switch-pc-game-controllers-23129289")
$ws.Hyperlinks.Add($ws.Range("B49"), "https://m.media-amazon.com/images/I/91j2tD6pM9L._SX342_.jpg")

# ---------------------------------------------------------------------------
# Selection follows the new row that was just filled in (matches the
# workbook's recorded selection after the edit).
# ---------------------------------------------------------------------------
$ws.Range("A50:D50").Select()
